# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#    using the same header/layout as the existing per-quarter fund-holding
#    sheets (e.g. "2021-Q4"), and fill it with the new quarter's data.
# 2. Update the "总计" (totals) summary sheet with a new leading row for
#    2022-Q1 and shift all the previous rows down by one.

function Set-TextCell($cell, $val) {
    # Force the cell to be stored as TEXT even when the value looks like a
    # number (e.g. "000179", "2.37") - matches the source data which keeps
    # these fund figures as plain strings rather than numeric values.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right before "总计"
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item(5)                              # "2021-Q4" - layout template
$totalSheetName = $wb.Worksheets.Item($wb.Worksheets.Count).Name # "总计" - last sheet (by name)
$insertBefore = $wb.Worksheets.Item($wb.Worksheets.Count)

$newSheet = $wb.Worksheets.Add($insertBefore)
$newSheet.Name = "2022-Q1"

# NOTE: sheet object references resolve by the index they were fetched at,
# not by a stable identity - after inserting a sheet before it, the "总计"
# sheet slides to a new index, so its handle must be re-fetched by name.
$totalSheet = $wb.Worksheets.Item($totalSheetName)

# Copy the header row + row-index column formatting from the template sheet
$template.Range("A1:H5").Copy()
$newSheet.Range("A1:H5").PasteSpecial(-4122)

# Header row
Set-TextCell $newSheet.Cells.Item(1, 2) "基金代码"
Set-TextCell $newSheet.Cells.Item(1, 3) "基金名称"
Set-TextCell $newSheet.Cells.Item(1, 4) "基金规模"
Set-TextCell $newSheet.Cells.Item(1, 5) "股票总仓位"
Set-TextCell $newSheet.Cells.Item(1, 6) "仓位占比"
Set-TextCell $newSheet.Cells.Item(1, 7) "持有市值(亿元)"
Set-TextCell $newSheet.Cells.Item(1, 8) "仓位排名"

$rows = @(
    @{ A=0; B="000179"; C="广发美国房地产指数QDII-人民币";  D="2.37"; E="92.38"; F="5.08"; G="0.1204"; H=2 },
    @{ A=1; B="000180"; C="广发美国房地产指数QDII - 美元"; D="2.37"; E="92.38"; F="5.08"; G="0.1204"; H=2 },
    @{ A=2; B="320017"; C="诺安全球收益不动产(QDII)";      D="0.29"; E="93.32"; F="8.66"; G="0.0251"; H=3 },
    @{ A=3; B="070031"; C="嘉实全球房地产(QDII)";          D="0.60"; E="95.08"; F="2.68"; G="0.0161"; H=8 }
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row.A
    Set-TextCell $newSheet.Cells.Item($r, 2) $row.B
    Set-TextCell $newSheet.Cells.Item($r, 3) $row.C
    Set-TextCell $newSheet.Cells.Item($r, 4) $row.D
    Set-TextCell $newSheet.Cells.Item($r, 5) $row.E
    Set-TextCell $newSheet.Cells.Item($r, 6) $row.F
    Set-TextCell $newSheet.Cells.Item($r, 7) $row.G
    $newSheet.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2. Update "总计" sheet: add 2022-Q1 as the new first data row, and
#    shift the rest down by one row.
# ------------------------------------------------------------------
$totalRows = @(
    @{ A=0; B="2022-Q1"; C=4; D=0.28 },
    @{ A=1; B="2021-Q4"; C=4; D=0.33 },
    @{ A=2; B="2021-Q3"; C=6; D=0.39 },
    @{ A=3; B="2021-Q2"; C=5; D=0.38 },
    @{ A=4; B="2021-Q1"; C=6; D=0.41 },
    @{ A=5; B="2020-Q4"; C=5; D=0.35 }
)

$r = 2
foreach ($row in $totalRows) {
    $totalSheet.Cells.Item($r, 1).Value = $row.A
    Set-TextCell $totalSheet.Cells.Item($r, 2) $row.B
    $totalSheet.Cells.Item($r, 3).Value = $row.C
    $totalSheet.Cells.Item($r, 4).Value = $row.D
    $r = $r + 1
}
